# recode: kuroi_tower & mitsui_tower
# Insert two new building rows ("kuroi_tower" and "mitsui_tower") above the
# existing "tsuno_building" row on the "items" sheet, pushing the rest of the
# table (tsuno_building .. temple, plus the trailing blank row) down by two
# rows, and mirror the existing skyscraper-row formula pattern into the new
# rows' computed graphics columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Insert two blank rows before the current row 15 ("tsuno_building"),
# shifting everything below it (and row styles) down by two.
$ws.Rows("15:16").Insert()

# --- Row 15: kuroi_tower -----------------------------------------------
$ws.Range("A15").Value = "kuroi_tower"
$ws.Range("B15").Value = "kuroi_tower"
$ws.Range("C15").Value = 38
$ws.Range("D15").Value = "1X1"
$ws.Range("E15").Value = "skyscraper"
$ws.Range("F15").Value = "NAME_KUROI_TOWER"
$ws.Range("G15").Value = 220
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 2000
$ws.Range("J15").Value = "0xFFFF"
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = "4 only"
$ws.Range("M15").Formula = "=VLOOKUP(L15,dropdowns!E:F,2,0)"
$ws.Range("N15").Value = 4
$ws.Range("O15").Value = 5
$ws.Range("P15").Value = "bitmask(HOUSE_FLAG_NOT_SLOPED,HOUSE_FLAG_PROTECTED)"
$ws.Range("Q15").Value = 24
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = "[PASS, 8],[MAIL, 4],[GOOD, 4]"
$ws.Range("T15").Formula = '=IF(NOT(D15="1X1"),"none",IF(E15="skyscraper",CONCATENATE(A15,"_c"),IF(E15="landmark",CONCATENATE(A15,"_k"),IF(E15="house",CONCATENATE(A15,"_h"),A15))))'
$ws.Range("U15").Formula = '=IF(D15="1X1","none",IF(E15="skyscraper",CONCATENATE(A15,"_c_north"),IF(E15="landmark",CONCATENATE(A15,"_k_north"),IF(E15="house",CONCATENATE(A15,"_h_north"),CONCATENATE(A15,"_north")))))'
$ws.Range("V15").Formula = '=IF(OR(D15="1X1",D15="2X1"),"none",IF(E15="skyscraper",CONCATENATE(A15,"_c_east"),IF(E15="landmark",CONCATENATE(A15,"_k_east"),CONCATENATE(A15,"_east"))))'
$ws.Range("W15").Formula = '=IF(OR(D15="1X1",D15="1X2"),"none",IF(E15="skyscraper",CONCATENATE(A15,"_c_west"),IF(E15="landmark",CONCATENATE(A15,"_k_west"),CONCATENATE(A15,"_west"))))'
$ws.Range("X15").Formula = '=IF(NOT(D15="2X2"),"none",IF(E15="skyscraper",CONCATENATE(A15,"_c_south"),IF(E15="landmark",CONCATENATE(A15,"_k_south"),CONCATENATE(A15,"_south"))))'
$ws.Range("Y15").Value = "none"

# --- Row 16: mitsui_tower ------------------------------------------------
$ws.Range("A16").Value = "mitsui_tower"
$ws.Range("B16").Value = "mitsui_tower"
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = "1X1"
$ws.Range("E16").Value = "skyscraper"
$ws.Range("F16").Value = "NAME_MITSUI_TOWER"
$ws.Range("G16").Value = 220
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = "0xFFFF"
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = "4 only"
$ws.Range("M16").Formula = "=VLOOKUP(L16,dropdowns!E:F,2,0)"
$ws.Range("N16").Value = 4
$ws.Range("O16").Value = 5
$ws.Range("P16").Value = "bitmask(HOUSE_FLAG_NOT_SLOPED,HOUSE_FLAG_PROTECTED)"
$ws.Range("Q16").Value = 24
$ws.Range("R16").Value = 10
$ws.Range("S16").Value = "[PASS, 8],[MAIL, 4],[GOOD, 4]"
$ws.Range("T16").Formula = '=IF(NOT(D16="1X1"),"none",IF(E16="skyscraper",CONCATENATE(A16,"_c"),IF(E16="landmark",CONCATENATE(A16,"_k"),IF(E16="house",CONCATENATE(A16,"_h"),A16))))'
$ws.Range("U16").Formula = '=IF(D16="1X1","none",IF(E16="skyscraper",CONCATENATE(A16,"_c_north"),IF(E16="landmark",CONCATENATE(A16,"_k_north"),IF(E16="house",CONCATENATE(A16,"_h_north"),CONCATENATE(A16,"_north")))))'
$ws.Range("V16").Formula = '=IF(OR(D16="1X1",D16="2X1"),"none",IF(E16="skyscraper",CONCATENATE(A16,"_c_east"),IF(E16="landmark",CONCATENATE(A16,"_k_east"),CONCATENATE(A16,"_east"))))'
$ws.Range("W16").Formula = '=IF(OR(D16="1X1",D16="1X2"),"none",IF(E16="skyscraper",CONCATENATE(A16,"_c_west"),IF(E16="landmark",CONCATENATE(A16,"_k_west"),CONCATENATE(A16,"_west"))))'
$ws.Range("X16").Formula = '=IF(NOT(D16="2X2"),"none",IF(E16="skyscraper",CONCATENATE(A16,"_c_south"),IF(E16="landmark",CONCATENATE(A16,"_k_south"),CONCATENATE(A16,"_south"))))'
$ws.Range("Y16").Value = "none"

# Update the visible selection to match the edited area (also settles the
# frozen-pane scroll position).
$ws.Range("A16").Select()
